# carjacking-by-month-yoy: roll the "through" date from 2022-07-31 to 2022-08-01
# and add the new August row (commit: "Add data for 2022-08-09").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The workbook/sheet name reflects the new "through" date.
$ws.Name = "Through 2022-08-01"

# 2. July no longer needs the "(through 07-31)" qualifier now that the month is complete.
$ws.Range("A8").Value = "July"

# 3. Two of the year-over-year counts were revised for the 2022 column.
$ws.Range("I3").Value = 140
$ws.Range("I8").Value = 167

# 4. Insert a new row under July for the (partial) August figures; this pushes the
#    existing "Total" row from row 9 down to row 10.
$ws.Rows.Item(9).Insert()

# New row inherits most formatting from the row above on insert, but not the
# thin-box border used by every month label in column A -- restore it explicitly
# and make sure font/alignment match the other header cells in the column.
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4160
$ws.Range("A9").Borders.LineStyle = 1
$ws.Range("A9").Value = "August (through 08-01)"

# August only has data starting in 2016, so B9 (2015) stays blank.
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 2

# 5. Total row (now row 10) is bumped by the new August counts.
$ws.Range("B10").Value = 162
$ws.Range("C10").Value = 303
$ws.Range("D10").Value = 467
$ws.Range("E10").Value = 427
$ws.Range("F10").Value = 306
$ws.Range("G10").Value = 626
$ws.Range("H10").Value = 915
$ws.Range("I10").Value = 974

# 6. Column A is widened slightly to fit the longer "August (through 08-01)" label.
$ws.Columns.Item(1).ColumnWidth = 21.75

Write-Output "carjacking-by-month-yoy: rolled to 2022-08-01 and appended August row"
